$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.762.04'
$ws.Range("E2").Value = '  -0.82%  '

$ws.Range("D3").Value = '1.623.70'
$ws.Range("E3").Value = '  -0.95%  '

$ws.Range("E5").Value = '  +0.04%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.5096'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +0.16%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  -0.16%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.06305'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -0.64%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '19.33'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -1.19%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.07773'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +0.20%  '

$ws.Range("D12").Value = '1.627.51'
$ws.Range("E12").Value = '  -0.76%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '4.207'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -1.57%  '

$ws.Range("D14").Value = '1.846.37'
$ws.Range("E14").Value = '  -1.08%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.5507'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +1.52%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '63.35'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -1.03%  '

$ws.Range("D17").Value = '0.0₅7480'
$ws.Range("E17").Value = '  -2.58%  '

$ws.Range("D18").Value = '25.766.07'
$ws.Range("E18").Value = '  -0.89%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -0.12%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '4.397'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -0.30%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '193.25'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -2.74%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '9.760'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -1.28%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '5.983'
$cell.Style = "Normal"

$ws.Range("E24").Value = '  -0.08%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '1.872'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -0.32%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '141.36'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("E27").Value = '  +5.00%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '15.49'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -0.83%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '6.694'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -1.69%  '

$ws.Range("E30").Value = '  +0.13%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '0.04854'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -0.96%  '

$ws.Range("E32").Value = '  -0.76%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '3.150'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -0.45%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.531'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +0.72%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '2.368'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -0.05%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.8915'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -1.38%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '2.534'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -2.07%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.5490'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +0.91%  '

$ws.Range("D39").Value = '1.112.02'
$ws.Range("E39").Value = '  -2.50%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.01540'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -1.35%  '

$ws.Range("E41").Value = '  -0.08%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '5.519'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +2.08%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.7954'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -1.70%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '97.03'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -2.12%  '

$ws.Range("D45").Value = '1.771.24'
$ws.Range("E45").Value = '  -0.35%  '

$ws.Range("D46").Value = '0.0₈115'
$ws.Range("E46").Value = '  -8.18%  '

$ws.Range("E47").Value = '  -2.36%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.9931'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -1.03%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '54.47'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -0.72%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.05121'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +0.10%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '7.537'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +2.92%  '

